$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and Column E (Volume 1h %) updates.
# Apostrophe-prefix is used on D-column values that would otherwise be
# auto-parsed by Excel as numbers (losing formatting like trailing zeros
# or turning tiny decimals into scientific notation); it forces them to
# stay plain text, matching the original inline-string cell content.

$ws.Range('D2').Value = '27.010.20'
$ws.Range('E2').Value = '  -2.61%  '
$ws.Range('D3').Value = '1.862.80'
$ws.Range('E3').Value = '  -2.36%  '
$ws.Range("D4").Value = "'0.9984"
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range("D5").Value = "'305.90"
$ws.Range('E5').Value = '  -2.08%  '
$ws.Range("D6").Value = "'0.9986"
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range("D7").Value = "'0.5066"
$ws.Range('E7').Value = '  -3.17%  '
$ws.Range('E8').Value = '  -1.11%  '
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range("D10").Value = "'20.63"
$ws.Range('E10').Value = '  -2.80%  '
$ws.Range("D11").Value = "'0.8826"
$ws.Range('E11').Value = '  -1.71%  '
$ws.Range("D12").Value = "'0.07569"
$ws.Range('E12').Value = '  -0.72%  '
$ws.Range('D13').Value = '1.855.69'
$ws.Range('E13').Value = '  -2.40%  '
$ws.Range("D14").Value = "'5.307"
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range("D15").Value = "'89.26"
$ws.Range('E15').Value = '  -3.15%  '
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range("D17").Value = "'0.000008406"
$ws.Range('E17').Value = '  -3.65%  '
$ws.Range("D18").Value = "'14.04"
$ws.Range('E18').Value = '  -2.92%  '
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('D20').Value = '27.033.88'
$ws.Range('E20').Value = '  -2.59%  '
$ws.Range('E21').Value = '  -2.03%  '
$ws.Range('D22').Value = '2.106.78'
$ws.Range('E22').Value = '  -0.84%  '
$ws.Range('E23').Value = '  -2.93%  '
$ws.Range("D24").Value = "'6.459"
$ws.Range('E24').Value = '  -1.91%  '
$ws.Range('E25').Value = '  -1.55%  '
$ws.Range("D26").Value = "'148.37"
$ws.Range('E26').Value = '  -3.33%  '
$ws.Range("D27").Value = "'17.97"
$ws.Range('E27').Value = '  -1.81%  '
$ws.Range("D28").Value = "'2.101"
$ws.Range('E28').Value = '  -2.85%  '
$ws.Range("D29").Value = "'112.77"
$ws.Range('E29').Value = '  -1.68%  '
$ws.Range("D30").Value = "'4.675"
$ws.Range('E30').Value = '  -3.54%  '
$ws.Range("D31").Value = "'4.713"
$ws.Range('E31').Value = '  -3.27%  '
$ws.Range("D32").Value = "'0.09046"
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range("D33").Value = "'0.05135"
$ws.Range('E33').Value = '  -2.71%  '
$ws.Range("D34").Value = "'3.038"
$ws.Range('E34').Value = '  -4.17%  '
$ws.Range("D35").Value = "'1.153"
$ws.Range('E35').Value = '  -6.80%  '
$ws.Range("D36").Value = "'0.7285"
$ws.Range('E36').Value = '  -5.37%  '
$ws.Range('E37').Value = '  -2.07%  '
$ws.Range("D38").Value = "'3.036"
$ws.Range('E38').Value = '  -1.05%  '
$ws.Range("D39").Value = "'2.469"
$ws.Range('E39').Value = '  -5.95%  '
$ws.Range("D40").Value = "'1.075"
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range("D41").Value = "'0.5291"
$ws.Range('E41').Value = '  -3.69%  '
$ws.Range("D42").Value = "'6.542"
$ws.Range('E42').Value = '  -1.57%  '
$ws.Range("D43").Value = "'115.84"
$ws.Range('E43').Value = '  +1.81%  '
$ws.Range("D44").Value = "'8.270"
$ws.Range('E44').Value = '  -2.55%  '
$ws.Range("D45").Value = "'0.1469"
$ws.Range('E45').Value = '  -2.57%  '
$ws.Range("D46").Value = "'0.9981"
$ws.Range('E46').Value = '  -0.13%  '
$ws.Range("D47").Value = "'0.4618"
$ws.Range('E47').Value = '  -3.53%  '
$ws.Range("D48").Value = "'9.985"
$ws.Range('E48').Value = '  -4.33%  '
$ws.Range('E49').Value = '  -3.20%  '
$ws.Range("D50").Value = "'36.53"
$ws.Range('E50').Value = '  -0.77%  '
$ws.Range("D51").Value = "'63.85"
$ws.Range('E51').Value = '  -3.94%  '
